# Regional Availability Factor.xlsx - "updated 4.0 files and mdl"
$wb = $excel.ActiveWorkbook

# --- About sheet: bump the "last updated" date from 2024-03-15 to 2024-03-28 ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = (Get-Date -Year 2024 -Month 3 -Day 28 -Hour 0 -Minute 0 -Second 0).Date

# --- RAF-capacity sheet: hydrogen combustion turbine / combined cycle capacity credit 0.3 -> 1 ---
$wsCapacity = $wb.Worksheets.Item("RAF-capacity")
$wsCapacity.Range("B24").Value = 1
$wsCapacity.Range("B25").Value = 1
$wsCapacity.Columns.Item(1).ColumnWidth = 29.04296875

# --- View state: RAF-capacity becomes the active/selected tab (was RAF-generation) ---
$wsCapacity.Select() | Out-Null
$wsCapacity.Range("B25").Select() | Out-Null
$excel.ActiveWindow.TopLeftCell = $wsCapacity.Range("A14")
$excel.ActiveWindow.Zoom = 80
